# Applies the commit "23/12/2017 MAMATHA CHICK IN":
#  1. Merge the two "THU Dec 14 ... 10:18:41 PST 2017" runs into one run.
#  2. Append a new "THU Dec 21 ... 11:30:40 PST 2017" purchase-details
#     block (mirroring the existing blocks) at the end of the document's
#     content, just after the "Amount balance - 21713.0" paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge "THU Dec 14" / " 10:18:41 PST 2017" into a single run.
#    A self find/replace over the (currently two-run) text coalesces it
#    into one run, matching the diff exactly.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("THU Dec 14 10:18:41 PST 2017", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "THU Dec 14 10:18:41 PST 2017", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Append the new "THU Dec 21" block after the "Amount balance" /
#    "- 21713.0" paragraph.
# ---------------------------------------------------------------------

# Locate the "Amount balance" paragraph that holds "- 21713.0" (the last
# purchase block currently in the document).
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*21713.0*") {
        $anchorIndex = $i
    }
}

function Add-EmptyParagraphAfter($paraIndex) {
    $p = $d.Paragraphs($paraIndex)
    $p.Range.InsertParagraphAfter()
    return ($paraIndex + 1)
}

function Add-TextRun($paraIndex, [string]$text) {
    $p = $d.Paragraphs($paraIndex)
    $pos = $p.Range.End - 1
    $r = $d.Range($pos, $pos)
    $r.InsertBefore($text)
}

function Add-Tab($paraIndex) {
    Add-TextRun $paraIndex "`t"
}

function Set-ParaBold($paraIndex, [bool]$bold) {
    $p = $d.Paragraphs($paraIndex)
    $p.Range.Font.Bold = $bold
}

# --- paragraph: empty, bold (blank line separating the previous block) ---
$idx = Add-EmptyParagraphAfter $anchorIndex
Set-ParaBold $idx $true

# --- paragraph: date/time stamp ---
$idx = Add-EmptyParagraphAfter $idx
Set-ParaBold $idx $false
Add-TextRun $idx "THU Dec 21"
Add-TextRun $idx " 11:30:40 PST 2017"

# --- paragraph: Person Name ---
$idx = Add-EmptyParagraphAfter $idx
Add-TextRun $idx "Person Name"
Add-Tab $idx
Add-Tab $idx
Add-Tab $idx
Add-TextRun $idx "`t- TRK"

# --- paragraph: Bill number ---
$idx = Add-EmptyParagraphAfter $idx
Add-TextRun $idx "Bill number"
Add-Tab $idx
Add-Tab $idx
Add-Tab $idx
Add-TextRun $idx "`t- 2069"

# --- paragraph: separator ---
$idx = Add-EmptyParagraphAfter $idx
Add-TextRun $idx "---------------------------------------------------------------"

# --- paragraph: Item Name ---
$idx = Add-EmptyParagraphAfter $idx
Add-TextRun $idx "Item Name"
Add-Tab $idx
Add-Tab $idx
Add-Tab $idx
Add-TextRun $idx "`t- CHOWCHOW"

# --- paragraph: Number of Pockets ---
$idx = Add-EmptyParagraphAfter $idx
Add-TextRun $idx "Number of Pockets"
Add-Tab $idx
Add-Tab $idx
Add-TextRun $idx "`t- 1"

# --- paragraph: Number of KGs ---
$idx = Add-EmptyParagraphAfter $idx
Add-TextRun $idx "Number of KGs"
Add-Tab $idx
Add-Tab $idx
Add-TextRun $idx "`t- 78"

# --- paragraph: Rate ---
$idx = Add-EmptyParagraphAfter $idx
Add-TextRun $idx "Rate"
Add-Tab $idx
Add-Tab $idx
Add-Tab $idx
Add-Tab $idx
Add-TextRun $idx "`t- 6"

# --- paragraph: Total Price ---
$idx = Add-EmptyParagraphAfter $idx
Add-TextRun $idx "Total Price"
Add-Tab $idx
Add-Tab $idx
Add-Tab $idx
Add-TextRun $idx "`t- 468.0"

# --- paragraph: Amount balance (bold) ---
$idx = Add-EmptyParagraphAfter $idx
Set-ParaBold $idx $true
Add-TextRun $idx "Amount balance"
Add-Tab $idx
Add-Tab $idx
Add-TextRun $idx "`t- 22181.0"

# --- paragraph: empty, non-bold ---
$idx = Add-EmptyParagraphAfter $idx
Set-ParaBold $idx $false

# --- paragraph: empty, bold ---
$idx = Add-EmptyParagraphAfter $idx
Set-ParaBold $idx $true

Write-Output "done; final paragraph count: $($d.Paragraphs.Count)"
